$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) column cells whose new values look like
# numbers, so Excel stores them as text (matching the original inlineStr
# cells) instead of auto-converting to Number and losing formatting like
# trailing zeros (e.g. "1.00" -> 1, "0.0000180" -> 0.000018).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "61.633.33"
$ws.Range("E2").Value = "  +2.38%  "

$ws.Range("D3").Value = "3.391.92"
$ws.Range("E3").Value = "  +3.52%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "575.26"
$ws.Range("E5").Value = "  +3.11%  "

$ws.Range("D6").Value = "138.04"
$ws.Range("E6").Value = "  +7.57%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "3.390.76"
$ws.Range("E8").Value = "  +3.44%  "

$ws.Range("D9").Value = "0.479"
$ws.Range("E9").Value = "  +1.06%  "

$ws.Range("D10").Value = "7.53"
$ws.Range("E10").Value = "  +3.09%  "

$ws.Range("E11").Value = "  +8.14%  "

$ws.Range("D12").Value = "0.395"
$ws.Range("E12").Value = "  +6.46%  "

$ws.Range("D13").Value = "3.968.02"
$ws.Range("E13").Value = "  +3.52%  "

$ws.Range("E14").Value = "  +1.51%  "

$ws.Range("D15").Value = "0.0000180"
$ws.Range("E15").Value = "  +7.40%  "

$ws.Range("D16").Value = "3.404.43"
$ws.Range("E16").Value = "  +4.11%  "

$ws.Range("D17").Value = "25.35"
$ws.Range("E17").Value = "  +4.25%  "

$ws.Range("D18").Value = "61.700.12"
$ws.Range("E18").Value = "  +2.13%  "

$ws.Range("D19").Value = "14.07"
$ws.Range("E19").Value = "  +6.11%  "

$ws.Range("D20").Value = "5.91"
$ws.Range("E20").Value = "  +4.57%  "

$ws.Range("D21").Value = "9.37"
$ws.Range("E21").Value = "  +4.17%  "

$ws.Range("D22").Value = "388.43"
$ws.Range("E22").Value = "  +10.29%  "

$ws.Range("D23").Value = "0.573"
$ws.Range("E23").Value = "  +3.91%  "

$ws.Range("D24").Value = "3.526.87"
$ws.Range("E24").Value = "  +3.54%  "

$ws.Range("D25").Value = "0.0000128"
$ws.Range("E25").Value = "  +16.90%  "

$ws.Range("E26").Value = "  +0.09%  "

$ws.Range("E27").Value = "  +2.44%  "

$ws.Range("E28").Value = "  +12.98%  "

$ws.Range("E29").Value = "  +5.35%  "

$ws.Range("D30").Value = "0.997"
$ws.Range("E30").Value = "  -0.21%  "

$ws.Range("D31").Value = "8.33"
$ws.Range("E31").Value = "  +6.53%  "

$ws.Range("D32").Value = "0.158"
$ws.Range("E32").Value = "  +5.98%  "

$ws.Range("D33").Value = "2.16"
$ws.Range("E33").Value = "  +3.20%  "

$ws.Range("E34").Value = "  -0.10%  "

$ws.Range("D35").Value = "3.422.68"
$ws.Range("E35").Value = "  +3.59%  "

$ws.Range("E36").Value = "  +3.67%  "

$ws.Range("D37").Value = "5.53"
$ws.Range("E37").Value = "  +4.15%  "

$ws.Range("D38").Value = "6.99"
$ws.Range("E38").Value = "  +2.78%  "

$ws.Range("E39").Value = "  +4.55%  "

$ws.Range("D40").Value = "162.07"
$ws.Range("E40").Value = "  +2.57%  "

$ws.Range("D41").Value = "0.0796"
$ws.Range("E41").Value = "  +6.04%  "

$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.06%  "

$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "1.74"
$ws.Range("E43").Value = "  +12.22%  "

$ws.Range("B44").Value = "ONDO"
$ws.Range("C44").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D44").Value = "1.23"
$ws.Range("E44").Value = "  +7.84%  "

$ws.Range("D45").Value = "0.772"
$ws.Range("E45").Value = "  +5.18%  "

$ws.Range("E46").Value = "  +2.44%  "

$ws.Range("D47").Value = "41.37"
$ws.Range("E47").Value = "  +1.09%  "

$ws.Range("D48").Value = "24.73"
$ws.Range("E48").Value = "  +10.69%  "

$ws.Range("D49").Value = "6.99"
$ws.Range("E49").Value = "  +5.10%  "

$ws.Range("E50").Value = "  +7.10%  "

$ws.Range("D51").Value = "2.383.57"
$ws.Range("E51").Value = "  +10.38%  "
